# The title paragraph ("Distributions croisées avec gjoint") has every run
# carrying an explicit bold (w:b / w:bCs) rPr override. That's redundant:
# the "Title" paragraph style already renders bold, so the commit strips
# the direct per-run formatting (visual result is unchanged - the title
# stays bold because the style itself is bold).
$d = $word.ActiveDocument

$titleXmlFragment = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">Distributions</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">croisées</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">avec</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">gjoint</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertXML($titleXmlFragment) | Out-Null

$d.Save()
